$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The product "NO-MIGRAIN 30 F.C. TAB" (previously on row 9) is removed from
# the report. All rows below it move up by one: what was row 10 (ONDALENZ)
# becomes row 9, what was row 11 (VISCERALGINE) becomes row 10, the totals
# row (previously row 12) becomes row 11 (with an updated total), and the
# footer row (previously row 13) becomes row 12.
# ---------------------------------------------------------------------------

# 1) Row 9: overwrite NO-MIGRAIN's data with what used to be row 10's data
#    (ONDALENZ 4 MG 5 ORODISPERSIBLE FILMS).
$ws.Range("B9").Value = "ONDALENZ 4 MG 5 ORODISPERSIBLE FILMS"
$ws.Range("H9").Value = "0:4"
$ws.Range("L9").Value = 32
$ws.Range("N9").Value = "0:0"

# 2) Row 10: overwrite ONDALENZ's data with what used to be row 11's data
#    (VISCERALGINE 10MG/5ML SYRUP 120 ML).
$ws.Range("B10").Value = "VISCERALGINE 10MG/5ML SYRUP 120 ML"
$ws.Range("H10").Value = "0:0"
$ws.Range("L10").Value = 35
$ws.Range("N10").Value = "1:0"

# 3) Row 11 used to hold VISCERALGINE's data (now duplicated into row 10).
#    Turn it into the totals row that used to be row 12: unmerge the old
#    product-row merges, blank out columns A-J, and merge K11:N11.
$ws.Range("B11:G11").UnMerge()
$ws.Range("H11:K11").UnMerge()
$ws.Range("L11:M11").UnMerge()

$ws.Range("A11:J11").Value = $null
$ws.Range("A11:J11").Style = "Normal"

$ws.Range("K11").Value = 229
$ws.Range("K11:N11").Style = $ws.Range("K12").Style
$ws.Range("K11:N11").Merge()
$ws.Rows(11).RowHeight = 25.5

# 4) Row 12 used to be the totals row (K12:N12 = 260). Clear it out; the
#    footer text that used to live on row 13 now belongs here instead.
$ws.Range("K12:N12").UnMerge()
$ws.Range("K12:N12").Value = $null
$ws.Range("K12:N12").Style = "Normal"

$ws.Range("A12").Value = "Saturday, 10 January, 2026 10:20 AM"
$ws.Range("A12").Style = $ws.Range("A13").Style
$ws.Range("F12").Value = "1/1"
$ws.Range("F12").Style = $ws.Range("F13").Style
$ws.Range("G12").Style = $ws.Range("G13").Style
$ws.Range("H12").Style = $ws.Range("H13").Style
$ws.Range("I12").Value = "developed by : Abdelaziz Talaat"
$ws.Range("I12").Style = $ws.Range("I13").Style
$ws.Range("J12:N12").Style = $ws.Range("J13:N13").Style
$ws.Rows(12).RowHeight = 16.5

$ws.Range("A13:E13").UnMerge()
$ws.Range("F13:G13").UnMerge()
$ws.Range("I13:N13").UnMerge()

$ws.Range("A12:E12").Merge()
$ws.Range("F12:G12").Merge()
$ws.Range("I12:N12").Merge()

# 5) Drop the now-empty old row 13 entirely (the data that lived there has
#    been copied up onto row 12 already).
$ws.Range("A13:N13").Value = $null
$ws.Range("A13:N13").Style = "Normal"
